# Refresh the "cryptos" price list (Price / Volume(1h) columns) with the
# latest scrape, matching the automated "Updated cryptos list ... with
# GitHub Actions" commit. Rows 19/20 also swap places (ShibaInu <-> 
# BitcoinCash) because the source ranking reordered them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character (U+2083) used in a couple of tiny price cells;
# built once and interpolated via "..." so it is never combined with `+`
# (string `+` against a numeric-looking left operand gets coerced to a numeric
# add in this runtime, which mangles the Unicode codepoint).
$sub3 = [char]0x2083

# Every "Price" (column D) cell is forced to Text ("@") before the write so
# a value like "231.87" or "1.00" is stored verbatim instead of being
# auto-parsed into a float (which would also strip the trailing zero /
# rewrite "1.00" as "1"). The style is then reset back to "Normal" so the
# cell doesn't keep a stray text-format style index it didn't have before.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.764.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.72%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.00%  '

# Row 4
$ws.Range("E4").Value = '  +0.39%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.88%  '

# Row 7
$ws.Range("E7").Value = '  +0.38%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.79'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.311'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.40%  '

# Row 10
$ws.Range("E10").Value = '  +3.45%  '

# Row 11
$ws.Range("E11").Value = '  +4.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.141.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.20%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.868.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.83%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.684'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.01%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.782.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.36%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '249.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.55%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0${sub3}0808"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.55%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.23%  '

# Row 23
$ws.Range("E23").Value = '  +0.32%  '

# Row 24
$ws.Range("E24").Value = '  +1.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.92%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.16%  '

# Row 28
$ws.Range("E28").Value = '  +2.31%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.98%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.370.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +38.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0552'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.61%  '

# Row 34
$ws.Range("E34").Value = '  +6.02%  '

# Row 35
$ws.Range("E35").Value = '  +4.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '98.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +20.50%  '

# Row 37
$ws.Range("E37").Value = '  +7.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.61%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.366.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.37%  '

# Row 41
$ws.Range("E41").Value = '  +6.03%  '

# Row 42
$ws.Range("E42").Value = '  +7.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.14%  '

# Row 44
$ws.Range("E44").Value = '  +3.95%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.42%  '

# Row 46
$ws.Range("E46").Value = '  +0.95%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.79%  '

# Row 48
$ws.Range("E48").Value = '  +2.40%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.039.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.21%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '105.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.88%  '

# Row 51
$ws.Range("E51").Value = '  +0.35%  '

